$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two changed threshold values (column C: ratio_threshold_range rows)
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 9

# Move the active selection from C7 (out of used range) to C3, matching the
# saved cursor position in the edited workbook.
$ws.Activate()
$ws.Range("C3").Select()
